$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.9202686041625157
$ws.Cells.Item(2, 3).Value = 0.2312851940432097
$ws.Cells.Item(2, 4).Value = 0.2216154181053156
$ws.Cells.Item(2, 5).Value = 0.1739950061645459
$ws.Cells.Item(2, 6).Value = 1.158484340464028
$ws.Cells.Item(2, 9).Value = 0.4743025055241681
$ws.Cells.Item(2, 10).Value = 0.181440027653359
$ws.Cells.Item(2, 14).Value = 1.002392765118174
$ws.Cells.Item(2, 15).Value = 2.593177007097211

$ws.Cells.Item(3, 2).Value = 0.8229987116440611
$ws.Cells.Item(3, 3).Value = 0.2035736324698121
$ws.Cells.Item(3, 4).Value = 0.2168581183873073
$ws.Cells.Item(3, 5).Value = 0.1700559920835012
$ws.Cells.Item(3, 6).Value = 1.152771200007365
$ws.Cells.Item(3, 9).Value = 0.4795121132207818
$ws.Cells.Item(3, 10).Value = 0.1771225843947235
$ws.Cells.Item(3, 14).Value = 1.006460705725516
$ws.Cells.Item(3, 15).Value = 2.588672475194073

$ws.Cells.Item(4, 2).Value = 0.7632777456278461
$ws.Cells.Item(4, 3).Value = 0.1865297834987416
$ws.Cells.Item(4, 4).Value = 0.2140201439043494
$ws.Cells.Item(4, 5).Value = 0.1677234901445317
$ws.Cells.Item(4, 6).Value = 1.150013388387926
$ws.Cells.Item(4, 9).Value = 0.4830496188176774
$ws.Cells.Item(4, 10).Value = 0.1745782619092537
$ws.Cells.Item(4, 14).Value = 1.009372486929486
$ws.Cells.Item(4, 15).Value = 2.587734933961968

$ws.Cells.Item(5, 2).Value = 0.7389432360144497
$ws.Cells.Item(5, 3).Value = 0.1795773785316896
$ws.Cells.Item(5, 4).Value = 0.21288461906002
$ws.Cells.Item(5, 5).Value = 0.1667946530482958
$ws.Cells.Item(5, 6).Value = 1.149078023406766
$ws.Cells.Item(5, 9).Value = 0.4845762209298066
$ws.Cells.Item(5, 10).Value = 0.1735682267089231
$ws.Cells.Item(5, 14).Value = 1.010663349112754
$ws.Cells.Item(5, 15).Value = 2.587811949816654

$ws.Cells.Item(6, 2).Value = 0.7349026898979218
$ws.Cells.Item(6, 3).Value = 0.1784225307842746
$ws.Cells.Item(6, 4).Value = 0.2126973355311748
$ws.Cells.Item(6, 5).Value = 0.1666417305556678
$ws.Cells.Item(6, 6).Value = 1.148934085324989
$ws.Cells.Item(6, 9).Value = 0.4848348438394758
$ws.Cells.Item(6, 10).Value = 0.1734021293204577
$ws.Cells.Item(6, 14).Value = 1.010883999665111
$ws.Cells.Item(6, 15).Value = 2.58785244577922

$ws.Cells.Item(7, 2).Value = 0.7629495507722766
$ws.Cells.Item(7, 3).Value = 0.1864360483944836
$ws.Cells.Item(7, 4).Value = 0.2140047447755933
$ws.Cells.Item(7, 5).Value = 0.1677108757105437
$ws.Cells.Item(7, 6).Value = 1.150000010836642
$ws.Cells.Item(7, 9).Value = 0.483069862985154
$ws.Cells.Item(7, 10).Value = 0.1745645317340276
$ws.Cells.Item(7, 14).Value = 1.009389473461688
$ws.Cells.Item(7, 15).Value = 2.587734114741636

$ws.Cells.Item(8, 2).Value = 0.8867301493752393
$ws.Cells.Item(8, 3).Value = 0.2217364363938543
$ws.Cells.Item(8, 4).Value = 0.2199579240351284
$ws.Cells.Item(8, 5).Value = 0.1726189891446879
$ws.Cells.Item(8, 6).Value = 1.156358676151896
$ws.Cells.Item(8, 9).Value = 0.4760283728502905
$ws.Cells.Item(8, 10).Value = 0.1799292395573246
$ws.Cells.Item(8, 14).Value = 1.003709555393627
$ws.Cells.Item(8, 15).Value = 2.591243982878581

$ws.Cells.Item(9, 2).Value = 1.129435424059977
$ws.Cells.Item(9, 3).Value = 0.290719381029561
$ws.Cells.Item(9, 4).Value = 0.2322874718304462
$ws.Cells.Item(9, 5).Value = 0.1829259138570833
$ws.Cells.Item(9, 6).Value = 1.174788671659812
$ws.Cells.Item(9, 9).Value = 0.4649149978424632
$ws.Cells.Item(9, 10).Value = 0.1912965171574541
$ws.Cells.Item(9, 14).Value = 0.995849766821209
$ws.Cells.Item(9, 15).Value = 2.61266752243327

$ws.Cells.Item(10, 2).Value = 1.307680207747183
$ws.Cells.Item(10, 3).Value = 0.3412427635976201
$ws.Cells.Item(10, 4).Value = 0.2417421073335788
$ws.Cells.Item(10, 5).Value = 0.1909142272174051
$ws.Cells.Item(10, 6).Value = 1.191979143481902
$ws.Cells.Item(10, 9).Value = 0.4584025110627081
$ws.Cells.Item(10, 10).Value = 0.2001673305942404
$ws.Cells.Item(10, 14).Value = 0.992065850044213
$ws.Cells.Item(10, 15).Value = 2.6373257986711

$ws.Cells.Item(11, 2).Value = 1.388742744772173
$ws.Cells.Item(11, 3).Value = 0.3641906446387679
$ws.Cells.Item(11, 4).Value = 0.2461286260977857
$ws.Cells.Item(11, 5).Value = 0.1946386679148375
$ws.Cells.Item(11, 6).Value = 1.200595791623698
$ws.Cells.Item(11, 9).Value = 0.4558007153940231
$ws.Cells.Item(11, 10).Value = 0.2043163326975872
$ws.Cells.Item(11, 14).Value = 0.9907751797901483
$ws.Cells.Item(11, 15).Value = 2.650492070846241

$ws.Cells.Item(12, 2).Value = 1.41943461133161
$ws.Cells.Item(12, 3).Value = 0.3728750208106817
$ws.Cells.Item(12, 4).Value = 0.2478019101156121
$ws.Cells.Item(12, 5).Value = 0.1960620151282626
$ws.Cells.Item(12, 6).Value = 1.203973474671187
$ws.Cells.Item(12, 9).Value = 0.4548675317125515
$ws.Cells.Item(12, 10).Value = 0.2059038214717788
$ws.Cells.Item(12, 14).Value = 0.9903482356932187
$ws.Cells.Item(12, 15).Value = 2.655758922750948

$ws.Cells.Item(13, 2).Value = 1.412824805813159
$ws.Cells.Item(13, 3).Value = 0.3710049353031195
$ws.Cells.Item(13, 4).Value = 0.2474409972297167
$ws.Cells.Item(13, 5).Value = 0.1957548949851073
$ws.Cells.Item(13, 6).Value = 1.203240923791597
$ws.Cells.Item(13, 9).Value = 0.4550661912411407
$ws.Cells.Item(13, 10).Value = 0.205561199913447
$ws.Cells.Item(13, 14).Value = 0.9904374388545847
$ws.Cells.Item(13, 15).Value = 2.654612098670412

$ws.Cells.Item(14, 2).Value = 1.391267887626498
$ws.Cells.Item(14, 3).Value = 0.3649052257149492
$ws.Cells.Item(14, 4).Value = 0.2462660441470632
$ws.Cells.Item(14, 5).Value = 0.1947555075551008
$ws.Cells.Item(14, 6).Value = 1.200871374868484
$ws.Cells.Item(14, 9).Value = 0.4557228974114302
$ws.Cells.Item(14, 10).Value = 0.2044466085699383
$ws.Cells.Item(14, 14).Value = 0.9907388168662123
$ws.Cells.Item(14, 15).Value = 2.650919740319239

$ws.Cells.Item(15, 2).Value = 1.378062985704105
$ws.Cells.Item(15, 3).Value = 0.3611682511832441
$ws.Cells.Item(15, 4).Value = 0.2455479379308798
$ws.Cells.Item(15, 5).Value = 0.1941450435771799
$ws.Cells.Item(15, 6).Value = 1.199434906690996
$ws.Cells.Item(15, 9).Value = 0.4561319340344276
$ws.Cells.Item(15, 10).Value = 0.2037660190389232
$ws.Cells.Item(15, 14).Value = 0.9909314648552652
$ws.Cells.Item(15, 15).Value = 2.648694691370082

$ws.Cells.Item(16, 2).Value = 1.302382002015577
$ws.Cells.Item(16, 3).Value = 0.339742318766298
$ws.Cells.Item(16, 4).Value = 0.2414571512298949
$ws.Cells.Item(16, 5).Value = 0.1906726447961233
$ws.Cells.Item(16, 6).Value = 1.191432071721692
$ws.Cells.Item(16, 9).Value = 0.4585798196456921
$ws.Cells.Item(16, 10).Value = 0.1998984710505312
$ws.Cells.Item(16, 14).Value = 0.9921588536034704
$ws.Cells.Item(16, 15).Value = 2.636504640428456

$ws.Cells.Item(17, 2).Value = 1.255947413905631
$ws.Cells.Item(17, 3).Value = 0.3265888439129583
$ws.Cells.Item(17, 4).Value = 0.2389694293160574
$ws.Cells.Item(17, 5).Value = 0.1885656000784977
$ws.Cells.Item(17, 6).Value = 1.186726761655919
$ws.Cells.Item(17, 9).Value = 0.4601740488354622
$ws.Cells.Item(17, 10).Value = 0.197554966648795
$ws.Cells.Item(17, 14).Value = 0.9930220282897011
$ws.Cells.Item(17, 15).Value = 2.629526214108324

$ws.Cells.Item(18, 2).Value = 1.229237500509328
$ws.Cells.Item(18, 3).Value = 0.3190199987077449
$ws.Cells.Item(18, 4).Value = 0.2375466164430833
$ws.Cells.Item(18, 5).Value = 0.1873622055931037
$ws.Cells.Item(18, 6).Value = 1.184095361538581
$ws.Cells.Item(18, 9).Value = 0.4611249535400326
$ws.Cells.Item(18, 10).Value = 0.196217739038417
$ws.Cells.Item(18, 14).Value = 0.9935590455538232
$ws.Cells.Item(18, 15).Value = 2.625695818934162

$ws.Cells.Item(19, 2).Value = 1.220193694057514
$ws.Cells.Item(19, 3).Value = 0.3164567604052309
$ws.Cells.Item(19, 4).Value = 0.237066263592439
$ws.Cells.Item(19, 5).Value = 0.1869562215739222
$ws.Cells.Item(19, 6).Value = 1.183217284921966
$ws.Cells.Item(19, 9).Value = 0.4614527382894202
$ws.Cells.Item(19, 10).Value = 0.1957668130865216
$ws.Cells.Item(19, 14).Value = 0.9937478373308863
$ws.Cells.Item(19, 15).Value = 2.624430389819338

$ws.Cells.Item(20, 2).Value = 1.260890670251626
$ws.Cells.Item(20, 3).Value = 0.327989399615177
$ws.Cells.Item(20, 4).Value = 0.2392334182006834
$ws.Cells.Item(20, 5).Value = 0.1887890168460089
$ws.Cells.Item(20, 6).Value = 1.187219889452763
$ws.Cells.Item(20, 9).Value = 0.4600008257882635
$ws.Cells.Item(20, 10).Value = 0.1978033297982051
$ws.Cells.Item(20, 14).Value = 0.9929259468281089
$ws.Cells.Item(20, 15).Value = 2.630250091395396

$ws.Cells.Item(21, 2).Value = 1.397599818925983
$ws.Cells.Item(21, 3).Value = 0.3666970103288349
$ws.Cells.Item(21, 4).Value = 0.2466108259054778
$ws.Cells.Item(21, 5).Value = 0.195048699801653
$ws.Cells.Item(21, 6).Value = 1.201564253133029
$ws.Cells.Item(21, 9).Value = 0.4555285925504897
$ws.Cells.Item(21, 10).Value = 0.2047735471842458
$ws.Cells.Item(21, 14).Value = 0.9906486184257659
$ws.Cells.Item(21, 15).Value = 2.65199664099265

$ws.Cells.Item(22, 2).Value = 1.486918681281793
$ws.Cells.Item(22, 3).Value = 0.3919625013272139
$ws.Cells.Item(22, 4).Value = 0.2515034802023308
$ws.Cells.Item(22, 5).Value = 0.1992154198121625
$ws.Cells.Item(22, 6).Value = 1.211608037269684
$ws.Cells.Item(22, 9).Value = 0.4529092587716335
$ws.Cells.Item(22, 10).Value = 0.2094243033141225
$ws.Cells.Item(22, 14).Value = 0.9895204448726247
$ws.Cells.Item(22, 15).Value = 2.667847928427165

$ws.Cells.Item(23, 2).Value = 1.439250657970888
$ws.Cells.Item(23, 3).Value = 0.378480905239087
$ws.Cells.Item(23, 4).Value = 0.248885706132512
$ws.Cells.Item(23, 5).Value = 0.1969846516559102
$ws.Cells.Item(23, 6).Value = 1.206186211408294
$ws.Cells.Item(23, 9).Value = 0.4542794134584902
$ws.Cells.Item(23, 10).Value = 0.2069333816606189
$ws.Cells.Item(23, 14).Value = 0.9900896537946124
$ws.Cells.Item(23, 15).Value = 2.659237598997009

$ws.Cells.Item(24, 2).Value = 1.258655867698565
$ws.Cells.Item(24, 3).Value = 0.3273562293415466
$ws.Cells.Item(24, 4).Value = 0.2391140457333165
$ws.Cells.Item(24, 5).Value = 0.1886879852870678
$ws.Cells.Item(24, 6).Value = 1.186996716676347
$ws.Cells.Item(24, 9).Value = 0.4600790329071422
$ws.Cells.Item(24, 10).Value = 0.1976910134091128
$ws.Cells.Item(24, 14).Value = 0.9929692582594782
$ws.Cells.Item(24, 15).Value = 2.629922260868199

$ws.Cells.Item(25, 2).Value = 1.063785789667406
$ws.Cells.Item(25, 3).Value = 0.2720846572618711
$ws.Cells.Item(25, 4).Value = 0.2288821935892713
$ws.Cells.Item(25, 5).Value = 0.180064593882399
$ws.Cells.Item(25, 6).Value = 1.169163066708407
$ws.Cells.Item(25, 9).Value = 0.4676319549144559
$ws.Cells.Item(25, 10).Value = 0.1881303834649088
$ws.Cells.Item(25, 14).Value = 0.9976259095982627
$ws.Cells.Item(25, 15).Value = 2.605309621831395
